# "Generate Report for Handoff" — add the newly-handed-off file
# cde05f60-9318-4a0d-86ef-ebb0b4ccad6e...md as a fresh row on every
# sheet (Overview, zh-cn, de-de) of the localization status report.

$wb = $excel.ActiveWorkbook

$fileNameLong  = "cde05f60-9318-4a0d-86ef-ebb0b4ccad6eoooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$pathAndName   = "e2e\cde05f60-9318-4a0d-86ef-ebb0b4ccad6eoooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$extension     = ".md"
$status        = "Ready for handoff"
$hoGenDate     = "2016-10-19 12:21:20"

$sourcePath    = "e2e"
$priority      = "ht"
$contentDup    = "False"
$toBeLocalized = "True"
$hasMetadata   = "False"

$zhHandoffFile = "cde05f60-9318-4a0d-86ef-ebb0b4ccad6eooooooooooooooooooooooooooooooooooooo.3db5fff0bb58257b0f0fb38b552bafad79da47f9.zh-cn.xlf"
$zhHandoffDate = "2016-10-19 12:21:09"
$deHandoffFile = "cde05f60-9318-4a0d-86ef-ebb0b4ccad6eooooooooooooooooooooooooooooooooooooo.3db5fff0bb58257b0f0fb38b552bafad79da47f9.de-de.xlf"
$deHandoffDate = "2016-10-19 12:21:20"
$handbackDate  = "0001-01-01 00:00:00"

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/440f2dd66aa258d5f142f32688075cb5c92842e7/e2e/"

function Fill-LangRow {
    param($ws, $handoffFile, $handoffDate)

    $t = $ws.ListObjects.Item(1)
    $row = $t.ListRows.Add()
    $r = $row.Range.Row

    $h = $ws.Hyperlinks.Add($ws.Cells.Item($r, 1), ($githubBase + $fileNameLong), "", "", $fileNameLong)

    $ws.Cells.Item($r, 2).Value = $extension
    $ws.Cells.Item($r, 3).Value = $status
    $ws.Cells.Item($r, 4).Value = $sourcePath
    $ws.Cells.Item($r, 5).Value = $priority
    $ws.Cells.Item($r, 6).Value = "'" + $contentDup
    $ws.Cells.Item($r, 7).Value = $handoffFile
    $ws.Cells.Item($r, 8).Value = $handoffDate
    $ws.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Cells.Item($r, 9).Value = "'"
    $ws.Cells.Item($r, 10).Value = "'"
    $ws.Cells.Item($r, 11).Value = $handbackDate
    $ws.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Cells.Item($r, 12).Value = "'"
    $ws.Cells.Item($r, 13).Value = "'" + $toBeLocalized
    $ws.Cells.Item($r, 14).Value = "'"
    $ws.Cells.Item($r, 15).Value = "'" + $hasMetadata
    $ws.Cells.Item($r, 16).Value = "'"
}

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
Fill-LangRow $wsZh $zhHandoffFile $zhHandoffDate

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
Fill-LangRow $wsDe $deHandoffFile $deHandoffDate

# --- Overview sheet ---
$wsOv = $wb.Worksheets.Item("Overview")
$tOv = $wsOv.ListObjects.Item(1)
$rowOv = $tOv.ListRows.Add()
$ro = $rowOv.Range.Row

$wsOv.Cells.Item($ro, 1).Value = $fileNameLong
$hOv = $wsOv.Hyperlinks.Add($wsOv.Cells.Item($ro, 2), ($githubBase + $fileNameLong), "", "", $pathAndName)
$wsOv.Cells.Item($ro, 3).Value = $extension
$wsOv.Cells.Item($ro, 4).Value = "'"
$wsOv.Cells.Item($ro, 5).Value = $status
$wsOv.Cells.Item($ro, 6).Value = $status
$wsOv.Cells.Item($ro, 7).Value = $hoGenDate
$wsOv.Cells.Item($ro, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
